$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Mark "Create Test Passed" (B) and "Read Test Passed" (C) as TRUE
# for all data rows (2 through 24). "Update Test Passed" (D) and
# "Delete Test Passed" (E) remain unchanged (FALSE).
$ws.Range("B2:C24").Value = $true
